$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new weekly column AI for date 29_06_2021
$ws.Range("AI1").Value = "29_06_2021"
$ws.Range("AI2").Value = 13
$ws.Range("AI3").Value = 19
$ws.Range("AI4").Value = 37
$ws.Range("AI5").Value = 57
$ws.Range("AI6").Value = 140
$ws.Range("AI7").Value = 301
$ws.Range("AI8").Value = 427
$ws.Range("AI9").Value = 578
$ws.Range("AI10").Value = 180
$ws.Range("AI11").Value = 13
$ws.Range("AI12").Formula = "=SUM(AI2:AI11)"

# Update the selection to match the author's saved view state
$null = $ws.Range("AK16").Select()
